# Update processed/analyzed market-object data: fill in previously-empty
# rows 74-78 with their now-known values, then append new daily rows
# 79-87 (2025-08-29 .. 2025-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in values for rows that were placeholders ('[]') ---
$ws.Range("B74").Value = "['BTCUSD.SPOT']"
$ws.Range("B75").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']"
$ws.Range("B76").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']"
$ws.Range("B77").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']"
$ws.Range("B78").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']"

# --- Append new rows 79-87 ---
$dates = @(
    "2025-08-29",
    "2025-08-30",
    "2025-08-31",
    "2025-09-01",
    "2025-09-02",
    "2025-09-03",
    "2025-09-04",
    "2025-09-05",
    "2025-09-06"
)

$objects = @(
    "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']",
    "['BTCUSD.SPOT']",
    "['BTCUSD.SPOT']",
    "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']",
    "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']",
    "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']",
    "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']",
    "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']",
    "[]"
)

$startRow = 79
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    # Force text storage (so the "2025-08-29"-style date strings aren't
    # auto-coerced into date serial numbers), then drop back to the
    # workbook's default "Normal" style so no stray per-cell formatting
    # is introduced (matches every other row in the column).
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dates[$i]
    $dateCell.Style = "Normal"

    $objCell = $ws.Cells.Item($row, 2)
    $objCell.NumberFormat = "@"
    $objCell.Value = $objects[$i]
    $objCell.Style = "Normal"
}
